$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "foo123"
$ws.Range("C4").Value = "foo123456"
$ws.Range("A5").Value = "boo123"
$ws.Range("C5").Value = "boo123456"

foreach ($addr in @("A4", "C4", "A5", "C5")) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

$ws.Range("C6").Select()
